$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Fly Spray "
$ws.Range("B2").Value = 1435453
$ws.Range("C2").Value = "E"
$ws.Range("D2").Value = 22
$ws.Range("E2").Value = 21
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 20
